# Insert a new weekly price record at row 49 of the "Feria Lagunitas de
# Puerto Montt - Espárragos" sheet, pushing the existing rows 49-70 down
# to 50-71 (sheet dimension grows from A1:R70 to A1:R71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 49; everything below shifts down
# by one (old row 49 -> 50, ..., old row 70 -> 71).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record's data.
$ws.Range("A49").Value = 4
$ws.Range("B49").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C49").Value = "Los Lagos"
$ws.Range("D49").Value = 45205
$ws.Range("E49").Value = 10
$ws.Range("F49").Value = 300000000
$ws.Range("G49").Value = "Espárragos"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 400
$ws.Range("K49").Value = 1900
$ws.Range("L49").Value = 2200
$ws.Range("M49").Value = 2050
$ws.Range("N49").Value = '$/kilo'
$ws.Range("O49").Value = "Provincia de Linares"
$ws.Range("P49").Value = 2050
$ws.Range("Q49").Value = 1
$ws.Range("R49").Value = "Hortaliza"
